# Reformat the guidance text in C2/C3 (add markdown-style emphasis, bullet
# markers, checkmarks) and refresh the dependent row heights + active
# selection, matching the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update cell C2 (row 2) text: Comunicação Abertura / Não cumprimentou corretamente ---
$ws.Range("C2").Value = '**Razão da Falha**  
	*Uso incorreto do modelo de cumprimento* (ex: *BraHello* quando há histórico de *Chatbot*).  
**Comportamento Correto**  
	*BraHello*:  
  	Sem *interação prévia* com *Chatbot*.  
	*BraAcceptTransfer*:  
	Quando há *conversa prévia* com *Chatbot*.  
**Comportamento Incorreto**  
	Não usar *modelo específico* ou usar *modelo inadequado* ao contexto.  
**Exemplos**  
	✅ **Correto**  
		*Cenário*: Cliente relata *depósito não creditado* após interagir com *Chatbot*.  
		*Ação do Agente*: Usa *IngAcceptTransfer* e menciona: *"Vi que seu contato é sobre o depósito"*.  
	❌ **Incorreto**  
		*Cenário*: Histórico de *Chatbot* visível.  
		*Ação do Agente*: Usa *BraHello* ou mensagem genérica (*"Como posso ajudar?"*).  
**Notas**  
	*Personalização* (ex: incluir nome do cliente) é permitida, mas não substitui o *modelo obrigatório*.  
	*E-mails*: Selecionar *modelo de e-mail* correspondente à consulta.  
	*Atraso > 60 segundos* na abertura é avaliado em **Questão 9 (Tempos de Espera)**.  '

# --- Update cell C3 (row 3) text: Compreensão do Problema / Não reconheceu todas as preocupações ---
$ws.Range("C3").Value = '**Razão da Falha:**
	Não identificar *todas as questões* do cliente (ex: ignorar uma *reclamação* ou *dúvida*).
**Comportamento Correto:**
	Reconhecer explicitamente *cada preocupação* (ex: *"Entendo que você tem 3 questões: aposta, suspensão e saque"*).
	Afirmar que irá *investigar* ou fazer *perguntas* que demonstrem *entendimento*.
**Comportamento Incorreto:**
	Focar apenas em *uma questão* e omitir outras mencionadas.
**Exemplos:**
	✅ **Correto:**
		*Cliente*: *"Preciso resolver uma aposta, minha conta está suspensa e quero saber do meu saque."*
		*Agente*: *"Vou verificar sua aposta, a suspensão e o status do saque."*
	❌ **Incorreto:**
		*Agente*: *"Vou checar sua aposta"* (ignora *suspensão* e *saque*).
**Notas:**
	*Perguntar* sobre *mensagens de erro* ou *confirmar detalhes* não é considerado *falha* (ex: *"Qual erro aparece?"*).'

# --- Adjust row heights to match the new wrapped-text extents ---
$ws.Rows.Item(2).RowHeight = 378.75
$ws.Rows.Item(3).RowHeight = 335.25

# --- Update the view: scroll down one row and move the active selection to C3 ---
$ws.Range("C3").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
